# Fix Training Data Issue (#48)
# The "Date" column (BF) on the active sheet stored the source file's own
# name ("4-26-2012-13") instead of an actual date value. Normalize every
# row's BF cell to the ISO-ish "2013-04-26" string (NBA stats were off by
# one day because of how they were originally scraped/shown).
#
# NOTE: a plain `Range.Value = "2013-04-26"` assignment gets auto-detected
# as a date literal and silently rewritten into a date serial number (with
# a new number-format style attached) instead of being kept as literal
# text. Routing the text through a formula cell (whose evaluated result is
# already a String, not a "user typed this" literal) and pasting only the
# value back avoids that smart-conversion, so the destination cells keep
# their original (unstyled) text-cell shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "4-26-2012-13"
$newDate = "2013-04-26"

# Out-of-the-way scratch cell used purely to "launder" the replacement
# string through a formula result so Excel won't re-interpret it as a date.
$scratch = $ws.Range("ZZ1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # column BF
    if ($cell.Text -eq $oldDate) {
        $scratch.Formula = '="' + $newDate + '"'
        $scratch.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

$scratch.ClearContents()
$excel.CutCopyMode = $false
